$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-5
# from 2023-10-22 (45221) to 2023-10-25 (45224)
$ws.Range("C2").Value = 45224
$ws.Range("C3").Value = 45224
$ws.Range("C4").Value = 45224
$ws.Range("C5").Value = 45224
